# "started migrating spaghetti into classes"
# Adds a new "Sheet2" after "Sheet1" containing the Timer1/PWM pulse-length
# scratch-work calculations, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet, placed immediately after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Column widths -----------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 34.498697916666664
$ws2.Columns.Item(2).ColumnWidth = 23.998697916666668
$ws2.Columns.Item(3).ColumnWidth = 26.498697916666668
$ws2.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws2.Columns.Item(5).ColumnWidth = 10.330729166666666

# --- Row 5 / Row 1 headers ----------------------------------------------
# (entered in this particular order so new shared-string ids land on the
# same indices the source workbook ends up with)
$ws2.Range("A5").Value = "Pulse Duration"
$ws2.Range("A1").Value = "Clock Period"
$ws2.Range("B1").Formula = "=1/16000000"
$ws2.Range("B5").Value = "Timer 1 Prescaler"
$ws2.Range("D5").Value = "Timer 1 Period"
$ws2.Range("E5").Value = "Pulse length"
$ws2.Range("C5").Value = "Timer 1 Comparator (actual value is x - 1)"

# --- Rows 6-25: data table (A=count, B=prescaler, C=comparator) --------
$dataRows = @(
  @(6, 1, 1, 2),
  @(7, 255, 1, 2),
  @(8, 1, 1, 65536),
  @(9, 255, 1, 65536),
  @(10, 1, 8, 2),
  @(11, 255, 8, 2),
  @(12, 1, 8, 65536),
  @(13, 255, 8, 65536),
  @(14, 1, 64, 2),
  @(15, 255, 64, 2),
  @(16, 1, 64, 65536),
  @(17, 255, 64, 65536),
  @(18, 1, 256, 2),
  @(19, 255, 256, 2),
  @(20, 1, 256, 65536),
  @(21, 255, 256, 65536),
  @(22, 1, 1024, 2),
  @(23, 255, 1024, 2),
  @(24, 1, 1024, 65536),
  @(25, 255, 1024, 65536)
)

foreach ($row in $dataRows) {
  $r = $row[0]
  $ws2.Cells.Item($r, 1).Value = $row[1]
  $ws2.Cells.Item($r, 2).Value = $row[2]
  $ws2.Cells.Item($r, 3).Value = $row[3]
  $ws2.Cells.Item($r, 4).Formula = "=2 * `$B`$1*B$r *C$r"
}

# Rows 22-25 got an explicit (custom) row height of 16 at some point.
$ws2.Range("A22:A25").RowHeight = 16

# --- Column E: "pulse length" = A * D, entered/filled in the same bursts
#     the original author used (reproduces the shared-formula grouping).
$ws2.Range("E7:E23").Formula = "=A7 *D7"
$ws2.Range("E12:E13").Formula = "=A12 *D12"
$ws2.Range("E16:E19").Formula = "=A16 *D16"
$ws2.Range("E24:E25").Formula = "=A24 *D24"
$ws2.Range("E6").Formula = "=A6 *D6"

# --- Row 28: sanity-check constant --------------------------------------
$ws2.Range("C28").Formula = "=2^16"

# --- Page setup (Excel defaults for a freshly inserted sheet) ----------
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36
$ws2.PageSetup.Orientation = 1

# --- Make Sheet2 the active tab / set its selection ---------------------
$ws2.Activate()
$ws2.Range("A6:A9").EntireRow.Select()
